# Append: 2025-10-27 12:51 JST
# Update the "取得日時" (retrieved-at) timestamp column (A2:A13) on the
# "ランサーズ" sheet from the previous run's timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-27 12:51:16"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
